$d = $word.ActiveDocument

$replacements = @(
    @{old="41×98=4018"; new="72×70=5040"},
    @{old="42×67=2814"; new="93×94=8742"},
    @{old="51×83=4233"; new="77×60=4620"},
    @{old="36×60=2160"; new="47×92=4324"},
    @{old="44×83=3652"; new="45×83=3735"},
    @{old="68×24=1632"; new="47×57=2679"},
    @{old="33×34=1122"; new="48×62=2976"},
    @{old="61×92=5612"; new="26×13=338"},
    @{old="72×57=4104"; new="72×72=5184"},
    @{old="34×25=850";  new="59×67=3953"},
    @{old="17×27=459";  new="58×82=4756"},
    @{old="69×92=6348"; new="41×37=1517"},
    @{old="40×54=2160"; new="85×24=2040"},
    @{old="26×21=546";  new="40×45=1800"},
    @{old="87×60=5220"; new="13×83=1079"},
    @{old="31×16=496";  new="79×62=4898"},
    @{old="60×83=4980"; new="87×51=4437"},
    @{old="90×58=5220"; new="89×49=4361"},
    @{old="62×69=4278"; new="38×97=3686"},
    @{old="38×58=2204"; new="14×54=756"},
    @{old="45×82=3690"; new="42×11=462"},
    @{old="56×55=3080"; new="47×13=611"},
    @{old="68×53=3604"; new="35×28=980"},
    @{old="66×63=4158"; new="46×62=2852"},
    @{old="18×85=1530"; new="55×68=3740"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
